# Remove the sample/demo rows that were used to reproduce the
# "login failed" error (Sakib Al Hasan, Mushfiq Rahim, Tahim Iqbal) from
# the "Bulk Upload" sample sheet, together with the mailto hyperlink that
# was attached to the sample email address, and move the active
# selection to E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bulk Upload")

# Clear the sample data (values + formatting) from the row cells that held
# the example names/emails/phone numbers/districts/lead sources, so the
# cells disappear from the sheet entirely, just like the original empty
# template rows below them.
$ws.Range("A2:B4").Clear()
$ws.Range("D2:F4").Clear()

# Column C (email) keeps its formatting/style, only the value (and its
# attached hyperlink) needs to go.
$ws.Range("C2").ClearContents()

# Remove the mailto: hyperlink that pointed at the sample email address.
$ws.Hyperlinks.Delete()

# Update the remembered selection/active cell for the sheet.
$ws.Range("E8").Select()
